$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# New row 13 (DRA051 / OPQA-4188||OPQA-4190) added beneath row 12.
# Duplicate row 12's formatting down into row 13 first (closest
# reachable match to its cell styles/border/wrap formatting), then
# overwrite the values with the new test-case content.
# ------------------------------------------------------------------
$ws.Range("A12:E12").Copy($ws.Range("A13:E13"))
$ws.Application.CutCopyMode = $false

$ws.Range("A13").Value = "DRA051"
$ws.Range("B13").Value = "OPQA-4188||OPQA-4190"
$ws.Range("C13").Value = "Verify that  if the STeAM account is pre-existing and missing a first or last name, the system should populate with a null value in the user's DRA profile.||Verify that a user who has a STeAM account with a missing first name, last name, or both should still be able to log into DRA."
$ws.Range("D13").Value = "Y"
# E13 stays blank, same as the copied-down E12 cell.

$ws.Rows.Item(13).RowHeight = 60

# New hyperlink on B13 pointing at the OPQA-4221 Jira issue (mirrors
# the hyperlink already present on B12), while keeping the cell's own
# displayed text as the OPQA-4188/4190 ids rather than the raw URL
# that Hyperlinks.Add would otherwise stamp into the cell.
$lnk = $ws.Hyperlinks.Add($ws.Range("B13"), "http://jira.bjz.apac.ime.reuters.com/browse/OPQA-4221")
$lnk.TextToDisplay = "http://jira.bjz.apac.ime.reuters.com/browse/OPQA-4221"

# Adding the hyperlink re-styles B13 with the built-in Hyperlink look;
# restore B12's plain cell formatting (same border/wrap, no hyperlink
# font) and put the intended text back afterwards.
$ws.Range("B12").Copy()
$ws.Range("B13").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false
$ws.Range("B13").Value = "OPQA-4188||OPQA-4190"

# Update the selection to reflect the newly added row (the file was
# scrolled down and the new row's data cell selected).
$ws.Range("C16:C18").Select()
